$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Removed rouge-2 and rouge-l: delete the ROUGE-2 (G) and ROUGE-L (H) columns entirely
$ws.Range("G1:H1").EntireColumn.Delete()

# Updated prompt: new Precision/Recall/F1/ROUGE-1 values for the data row
$ws.Range("C2").Value = 0.7666666666666666
$ws.Range("D2").Value = 0.9199999999999999
$ws.Range("E2").Value = 0.8033333333333333
$ws.Range("F2").Value = 0.7530091485447581
